$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 15307.85069382231
$ws.Range("D3").Value = 833.5794467332574
$ws.Range("E3").Value = 2825.673217591321

$ws.Range("B4").Value = 6804.466629161197
$ws.Range("D4").Value = 224.4594209241195
$ws.Range("E4").Value = 1361.059203808987

$ws.Range("B5").Value = 5112.259786301379
$ws.Range("D5").Value = 234.5908054794514
$ws.Range("E5").Value = 663.2953342465759

$ws.Range("B6").Value = 9399.504675342459
$ws.Range("D6").Value = 709.309699315068
$ws.Range("E6").Value = 1428.166626712329

$ws.Range("B7").Value = 14569.99556164384
$ws.Range("D7").Value = 882.369920547945
$ws.Range("E7").Value = 2705.965230136986

$ws.Range("B8").Value = 19731.90009657538
$ws.Range("D8").Value = 972.9928547945178
$ws.Range("E8").Value = 4314.070471917806

$ws.Range("B9").Value = 28230.24325753423
$ws.Range("D9").Value = 1233.665660273973
$ws.Range("E9").Value = 5027.55367123287

$ws.Range("F10").Value = 22043304.99910413

$ws.Range("G11").Value = 0.7609558168867352

$ws.Range("F12").Value = 1200354.403295889
$ws.Range("G12").Value = 0.05445437530101195

$ws.Range("F13").Value = 4068969.433331504
$ws.Range("G13").Value = 0.1845898078122529
